{"js": "// Commit: \"faltaba sacar una cosa\" (a leftover sentence needed to be removed)\n//\n// 1) In the \"hotel_usuario_rol ... con 3 FKs.\" paragraph, the trailing\n//    sentence \" Por cuestiones de tiempo, se decidi\u00f3 realizar que la\n//    implementaci\u00f3n del rol sea independiente del hotel. Es decir, un\n//    mismo usuario tendr\u00e1 los mismos roles en todos sus hoteles\n//    asignados. \" is deleted, leaving a single trailing space after\n//    \"con 3 FKs.\".\n// 2) In the \"Para realizar el alta de un usuario nuevo...\" paragraph,\n//    the two literal tab characters (before each \"Funcionalidades -> ...\")\n//    are replaced with plain spaces.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Locate the two affected paragraphs by their (stable) content instead of\n// a hard-coded index, so the script keeps working even if earlier\n// paragraphs/images shift the indices around.\nlet fksParagraph = null;\nlet altaParagraph = null;\nfor (let i = 0; i < paragraphs.items.length; i++) {\n  const text = paragraphs.items[i].text;\n  if (fksParagraph === null && text.indexOf(\"hotel_usuario_rol\") !== -1) {\n    fksParagraph = paragraphs.items[i];\n  }\n  if (altaParagraph === null && text.indexOf(\"Para realizar el alta de un usuario nuevo\") !== -1) {\n    altaParagraph = paragraphs.items[i];\n  }\n}\n\n// --- 1) Remove the leftover sentence after \"con 3 FKs.\" -----------------\nif (fksParagraph) {\n  const staleSentence = fksParagraph.search(\n    \" Por cuestiones de tiempo, se decidi\u00f3 realizar que la implementaci\u00f3n del rol sea independiente del hotel. Es decir, un mismo usuario tendr\u00e1 los mismos roles en todos sus hoteles asignados. \",\n    { matchCase: true }\n  );\n  await context.sync();\n\n  if (staleSentence.items.length > 0) {\n    staleSentence.items[0].insertText(\" \", Word.InsertLocation.replace);\n    await context.sync();\n  }\n}\n\n// --- 2) Replace the two tabs with spaces ---------------------------------\nif (altaParagraph) {\n  const tabAfterIngresar = altaParagraph.search(\"ingresar en:\\t\", { matchCase: true });\n  const tabAfterUsuario = altaParagraph.search(\"usuario en:\\t\", { matchCase: true });\n  await context.sync();\n\n  if (tabAfterIngresar.items.length > 0) {\n    tabAfterIngresar.items[0].insertText(\"ingresar en: \", Word.InsertLocation.replace);\n  }\n  if (tabAfterUsuario.items.length > 0) {\n    tabAfterUsuario.items[0].insertText(\"usuario en: \", Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Commit: \"faltaba sacar una cosa\" (a leftover sentence needed to be removed)\n#\n# 1) In the \"hotel_usuario_rol ... con 3 FKs.\" paragraph, the trailing\n#    sentence \" Por cuestiones de tiempo, se decidi\u00f3 realizar que la\n#    implementaci\u00f3n del rol sea independiente del hotel. Es decir, un\n#    mismo usuario tendr\u00e1 los mismos roles en todos sus hoteles\n#    asignados. \" is deleted, leaving a single trailing space after\n#    \"con 3 FKs.\".\n# 2) In the \"Para realizar el alta de un usuario nuevo...\" paragraph,\n#    the two literal tab characters are replaced with plain spaces.\n\n$d = $word.ActiveDocument\n\n# --- 1) Remove the leftover sentence after \"con 3 FKs.\" -----------------\n$targetPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -match \"hotel_usuario_rol\") {\n        $targetPara = $d.Paragraphs.Item($i)\n        break\n    }\n}\n\n$find = $targetPara.Range.Find\n$find.ClearFormatting()\n$find.Replacement.ClearFormatting()\n$find.Text = \" Por cuestiones de tiempo, se decidi\u00f3 realizar que la implementaci\u00f3n del rol sea independiente del hotel. Es decir, un mismo usuario tendr\u00e1 los mismos roles en todos sus hoteles asignados. \"\n$find.Replacement.Text = \" \"\n$find.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n# --- 2) Replace the two tabs with spaces ---------------------------------\n$altaPara = $null\nfor ($i = 1; $i -le $d.Paragraphs.Count; $i++) {\n    if ($d.Paragraphs.Item($i).Range.Text -match \"Para realizar el alta\") {\n        $altaPara = $d.Paragraphs.Item($i)\n        break\n    }\n}\n\n$find1 = $altaPara.Range.Find\n$find1.ClearFormatting()\n$find1.Replacement.ClearFormatting()\n$find1.Text = \"ingresar en:^t\"\n$find1.Replacement.Text = \"ingresar en: \"\n$find1.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n\n$find2 = $altaPara.Range.Find\n$find2.ClearFormatting()\n$find2.Replacement.ClearFormatting()\n$find2.Text = \"usuario en:^t\"\n$find2.Replacement.Text = \"usuario en: \"\n$find2.Execute($null, $null, $null, $null, $null, $null, $null, $null, $null, $null, 2) | Out-Null\n"}
